$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (Dropout = 0.5): add Val Accuracy "0.4799" (entered/stored as text)
# and Val Loss 1.211
$ws.Range("B5").Formula = '="0.4799"'
$ws.Range("B5").Copy()
$ws.Range("B5").PasteSpecial(-4163)
$ws.Range("C5").Value = 1.2110000000000001

# Row 6 (Dropout = 0.75): add Val Accuracy "0.6066" (entered/stored as text)
# and Val Loss 0.6943
$ws.Range("B6").Formula = '="0.6066"'
$ws.Range("B6").Copy()
$ws.Range("B6").PasteSpecial(-4163)
$ws.Range("C6").Value = 0.69430000000000003

$excel.CutCopyMode = 0

# Move the active cell selection to F10
[void]$ws.Range("F10").Select()
